$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the move list ("Tackle,SwordsDance") for Bulbasuar (row 2, column L)
$ws.Range("L2").Value = "Tackle,SwordsDance"

# New content is left-aligned
$ws.Range("L2").HorizontalAlignment = -4131  # xlLeft

# Update current selection to D7
$ws.Range("D7").Select()
